$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '44.458.36'
$c.Style = "Normal"

$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  +2.73%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.367.94'
$c.Style = "Normal"

$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '0.675'
$c.Style = "Normal"

$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  +3.92%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '239.13'
$c.Style = "Normal"

$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +2.88%  '
$c.Style = "Normal"

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '73.74'
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +7.99%  '
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.549'
$c.Style = "Normal"

$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +19.98%  '
$c.Style = "Normal"

$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +7.15%  '
$c.Style = "Normal"

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '30.15'
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  +13.94%  '
$c.Style = "Normal"

$c = $ws.Range('B12')
$c.NumberFormat = "@"
$c.Value = 'TRON'
$c.Style = "Normal"

$c = $ws.Range('C12')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c.Style = "Normal"

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.107'
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +2.02%  '
$c.Style = "Normal"

$c = $ws.Range('B13')
$c.NumberFormat = "@"
$c.Value = 'WrappedliquidstakedEther2.0'
$c.Style = "Normal"

$c = $ws.Range('C13')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '2.716.40'
$c.Style = "Normal"

$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '16.83'
$c.Style = "Normal"

$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +7.86%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '6.80'
$c.Style = "Normal"

$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  +9.02%  '
$c.Style = "Normal"

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.905'
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +7.56%  '
$c.Style = "Normal"

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.361.24'
$c.Style = "Normal"

$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '44.488.93'
$c.Style = "Normal"

$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  +2.76%  '
$c.Style = "Normal"

$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +5.13%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '77.47'
$c.Style = "Normal"

$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +4.68%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.47'
$c.Style = "Normal"

$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +3.48%  '
$c.Style = "Normal"

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '254.38'
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  +2.39%  '
$c.Style = "Normal"

$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  -3.64%  '
$c.Style = "Normal"

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.52'
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  +2.87%  '
$c.Style = "Normal"

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '10.41'
$c.Style = "Normal"

$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +4.72%  '
$c.Style = "Normal"

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"

$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +1.22%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '22.49'
$c.Style = "Normal"

$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +0.85%  '
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +4.73%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '173.99'
$c.Style = "Normal"

$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"

$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  +2.44%  '
$c.Style = "Normal"

$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  +5.26%  '
$c.Style = "Normal"

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.0743'
$c.Style = "Normal"

$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +7.54%  '
$c.Style = "Normal"

$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +3.69%  '
$c.Style = "Normal"

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '5.20'
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +2.85%  '
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +7.16%  '
$c.Style = "Normal"

$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -3.18%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '6.49'
$c.Style = "Normal"

$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.Style = "Normal"

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.0273'
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  +7.22%  '
$c.Style = "Normal"

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '19.66'
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +8.70%  '
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '8.84'
$c.Style = "Normal"

$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -1.18%  '
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  +3.13%  '
$c.Style = "Normal"

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.0988'
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +4.20%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.Style = "Normal"

$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  +1.71%  '
$c.Style = "Normal"

$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  +12.83%  '
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  +0.43%  '
$c.Style = "Normal"

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '98.90'
$c.Style = "Normal"

$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '2.36'
$c.Style = "Normal"

$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  +4.41%  '
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '1.446.46'
$c.Style = "Normal"

$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '2.590.32'
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.Style = "Normal"

